$d = $word.ActiveDocument

# --- Tables: tighten the left indent / left cell margins a touch ---
# tblInd 44 -> 42 dxa (2.2pt -> 2.1pt), tblCellMar/tcMar left 38 -> 35 dxa (1.9pt -> 1.75pt)
foreach ($t in $d.Tables) {
    $t.Rows.LeftIndent = 2.1
    $t.LeftPadding = 1.75
    foreach ($r in $t.Rows) {
        foreach ($c in $r.Cells) {
            $c.LeftPadding = 1.75
        }
    }
}

# --- Normal style: stop hanging/overflow punctuation ---
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false

# --- New character styles ListLabel156 .. ListLabel164 (merged demonBehavior lists) ---
for ($i = 156; $i -le 164; $i++) {
    $ls = $d.Styles.Add("ListLabel " + $i, 2)
    $ls.QuickStyle = $true
    $ls.Font.NameBi = "OpenSymbol"
}
